$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = 'Shenzhen Forever Industrial Co., Ltd.'
$ws.Range('D3').Value = 'BEIJING DAORI PLASTICS CO., LTD.'
$ws.Range('D4').Value = 'Dongguan Taihong Packaging Co., Ltd.'
$ws.Range('D5').Value = 'Dongguan Ruitai Packaging Technology Co., Ltd.'
$ws.Range('D6').Value = 'Xinyu Yushui High-Tech Co., Ltd.'
$ws.Range('D7').Value = 'Paper Sailing Factory Company for Paper Products'
$ws.Range('D8').Value = 'Shenzhen Kingcolor Paper Co., Ltd.'
$ws.Range('D9').Value = 'Dongguan Kaijing New Material Technology Co., Ltd.'
$ws.Range('D10').Value = 'Huizhou Yi To Culture Media Co., Ltd.'
$ws.Range('D11').Value = 'Best Zaza (xiamen) Industrial & Trade Co., Ltd.'
$ws.Range('D12').Value = 'Zhuhai Topa Paper Co., Ltd.'
$ws.Range('D13').Value = 'Shandong Luxu Supply Chain Management Co., Ltd.'
$ws.Range('D14').Value = 'HOANG MINH EXPORT COMPANY LIMITED'
$ws.Range('D15').Value = 'Shenzhen Greentree Co., Ltd.'
$ws.Range('D16').Value = 'Dongguan Sunkey Paper&Printing Co., Ltd.'
$ws.Range('D17').Value = 'Fy Packing Material (shanghai) Co., Ltd.'
$ws.Range('D18').Value = 'Ningbo Kunpeng printing co.ltd.'
$ws.Range('D19').Value = 'Shenzhen Dayu Environmental Protection Packaging Technology Co., Ltd.'
$ws.Range('D20').Value = 'Shenzhen Sailing Paper Co., Ltd.'
$ws.Range('D21').Value = 'Shenzhen Deyun Printing Products Co., Ltd.'
$ws.Range('B22').Value = '0,0259 €'
$ws.Range('D22').Value = 'Anhui Jimei Digital Technology Co., Ltd.'
$ws.Range('D23').Value = 'Shenzhen Forever Industrial Co., Ltd.'
$ws.Range('D24').Value = 'Shenzhen Kaiyu Electronic Technology Co., Ltd.'
$ws.Range('D25').Value = 'Qingdao Zhongcai Packaging Technology Co., Ltd.'
$ws.Range('D26').Value = 'Jiangsu OPT Barcode Label Co., Ltd.'
$ws.Range('D27').Value = 'Shenzhen Rongzhen Label Consumabels Corp., Ltd.'
$ws.Range('D28').Value = 'Jiangmen MST Packaging Co., Ltd.'
$ws.Range('D29').Value = 'Shenzhen Dinghao Paper Product Packaging Corp. Ltd.'
$ws.Range('D30').Value = 'Guangdong Juncheng Printing Technology Co., Ltd.'
$ws.Range('D31').Value = 'Shenzhen Fulida Printing Co., Ltd.'
$ws.Range('D32').Value = 'Shanghai Brightpac Printing Co., Ltd.'
$ws.Range('D33').Value = 'Shenzhen Xinxincai Security Label Technology Co., Ltd.'
$ws.Range('D34').Value = 'Shenzhen Ownlikes Technology Co., Ltd.'
$ws.Range('D35').Value = 'HONG KONG RUNYIN TONGDA CO., LIMITED'
$ws.Range('D36').Value = 'Shenzhen Shengtianxin Printing Co., Ltd.'
$ws.Range('D37').Value = 'Dongguan Caihe Printing Co., Ltd.'
$ws.Range('D38').Value = 'Guangdong Baobang Environmental Protection High-Tech Co., Ltd.'
$ws.Range('D39').Value = 'Shenzhen Enjoy Technology Co., Ltd.'
$ws.Range('D40').Value = 'Fujian Hongye Technology Co., Ltd.'
$ws.Range('D41').Value = 'Shenzhen Sailing Paper Company Limited'
$ws.Range('D42').Value = 'Hubei Bisheng Paper Industry Co., Ltd.'
$ws.Range('D43').Value = 'Shenzhen Coolmate Printing Co., Ltd.'
$ws.Range('D44').Value = 'Shenzhen He Li Packaging And Printing Products Co., Ltd.'
$ws.Range('D45').Value = 'Shenzhen Hangte Technology Development Co., Ltd.'
$ws.Range('D46').Value = 'HONGKONG PURE TECHNOLOGY LIMITED'
$ws.Range('B47').Value = '0,4317 €'
$ws.Range('D47').Value = 'Jinya New Materials Co., Ltd.'
$ws.Range('D48').Value = 'Dongguan Fenglin Printing Co., Ltd.'
$ws.Range('D49').Value = 'Zhejiang Jingran Trading Co., Limited'
